$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.171.53'
$ws.Range("E2").Value = '  -1.82%  '

$ws.Range("D3").Value = '3.134.81'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'594.81"
$ws.Range("E5").Value = '  -2.32%  '

$ws.Range("D6").Value = "'136.74"
$ws.Range("E6").Value = '  -4.45%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = '3.125.11'
$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("E9").Value = '  -2.25%  '

$ws.Range("E10").Value = '  -2.81%  '

$ws.Range("D11").Value = "'5.23"
$ws.Range("E11").Value = '  -2.38%  '

$ws.Range("E12").Value = '  -2.25%  '

$ws.Range("E13").Value = '  -2.47%  '

$ws.Range("D14").Value = "'34.30"
$ws.Range("E14").Value = '  -2.73%  '

$ws.Range("D15").Value = '3.640.60'
$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("E16").Value = '  +3.06%  '

$ws.Range("D17").Value = '63.160.61'
$ws.Range("E17").Value = '  -1.75%  '

$ws.Range("D18").Value = '3.128.79'
$ws.Range("E18").Value = '  -0.21%  '

$ws.Range("E19").Value = '  -1.85%  '

$ws.Range("D20").Value = "'474.27"
$ws.Range("E20").Value = '  -0.23%  '

$ws.Range("D21").Value = "'14.27"
$ws.Range("E21").Value = '  -3.36%  '

$ws.Range("E22").Value = '  -2.03%  '

$ws.Range("D23").Value = "'7.74"
$ws.Range("E23").Value = '  -0.41%  '

$ws.Range("D24").Value = "'86.63"
$ws.Range("E24").Value = '  +0.95%  '

$ws.Range("D25").Value = "'13.01"
$ws.Range("E25").Value = '  -3.63%  '

$ws.Range("E26").Value = '  -0.07%  '

$ws.Range("E27").Value = '  -1.37%  '

$ws.Range("D28").Value = "'7.15"
$ws.Range("E28").Value = '  -2.53%  '

$ws.Range("E29").Value = '  -5.62%  '

$ws.Range("E30").Value = '  +0.25%  '

$ws.Range("E31").Value = '  +0.05%  '

$ws.Range("D32").Value = "'26.85"
$ws.Range("E32").Value = '  +1.01%  '

$ws.Range("D33").Value = "'0.108"
$ws.Range("E33").Value = '  -5.37%  '

$ws.Range("E34").Value = '  -3.83%  '

$ws.Range("E35").Value = '  -2.02%  '

$ws.Range("E36").Value = '  -2.01%  '

$ws.Range("D37").Value = "'52.17"
$ws.Range("E37").Value = '  -0.78%  '

$ws.Range("D38").Value = '0.0₃0709'
$ws.Range("E38").Value = '  -4.26%  '

$ws.Range("E39").Value = '  -1.27%  '

$ws.Range("D40").Value = "'423.20"
$ws.Range("E40").Value = '  -5.56%  '

$ws.Range("D41").Value = "'8.25"
$ws.Range("E41").Value = '  -0.58%  '

$ws.Range("E42").Value = '  -9.38%  '

$ws.Range("D43").Value = '2.894.78'
$ws.Range("E43").Value = '  +0.60%  '

$ws.Range("E44").Value = '  -3.56%  '

$ws.Range("D45").Value = "'0.263"
$ws.Range("E45").Value = '  +0.58%  '

$ws.Range("D46").Value = "'2.13"
$ws.Range("E46").Value = '  -3.84%  '

$ws.Range("D48").Value = "'25.75"
$ws.Range("E48").Value = '  -2.08%  '

$ws.Range("E49").Value = '  -5.01%  '

$ws.Range("E50").Value = '  -0.59%  '

$ws.Range("D51").Value = "'120.37"
$ws.Range("E51").Value = '  -0.14%  '
